$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-05 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-06 Saturday", 2) | Out-Null
$d.Content.Find.Execute("44×36=1584", $true, $false, $false, $false, $false, $true, 1, $false, "71×41=2911", 2) | Out-Null
$d.Content.Find.Execute("42×72=3024", $true, $false, $false, $false, $false, $true, 1, $false, "36×74=2664", 2) | Out-Null
$d.Content.Find.Execute("23×20=460", $true, $false, $false, $false, $false, $true, 1, $false, "74×42=3108", 2) | Out-Null
$d.Content.Find.Execute("28×83=2324", $true, $false, $false, $false, $false, $true, 1, $false, "59×79=4661", 2) | Out-Null
$d.Content.Find.Execute("79×68=5372", $true, $false, $false, $false, $false, $true, 1, $false, "91×72=6552", 2) | Out-Null
$d.Content.Find.Execute("88×18=1584", $true, $false, $false, $false, $false, $true, 1, $false, "31×70=2170", 2) | Out-Null
$d.Content.Find.Execute("33×88=2904", $true, $false, $false, $false, $false, $true, 1, $false, "96×33=3168", 2) | Out-Null
$d.Content.Find.Execute("54×97=5238", $true, $false, $false, $false, $false, $true, 1, $false, "57×92=5244", 2) | Out-Null
$d.Content.Find.Execute("14×21=294", $true, $false, $false, $false, $false, $true, 1, $false, "44×93=4092", 2) | Out-Null
$d.Content.Find.Execute("50×91=4550", $true, $false, $false, $false, $false, $true, 1, $false, "25×13=325", 2) | Out-Null
$d.Content.Find.Execute("97×25=2425", $true, $false, $false, $false, $false, $true, 1, $false, "65×89=5785", 2) | Out-Null
$d.Content.Find.Execute("76×98=7448", $true, $false, $false, $false, $false, $true, 1, $false, "88×47=4136", 2) | Out-Null
$d.Content.Find.Execute("65×46=2990", $true, $false, $false, $false, $false, $true, 1, $false, "86×85=7310", 2) | Out-Null
$d.Content.Find.Execute("89×26=2314", $true, $false, $false, $false, $false, $true, 1, $false, "89×17=1513", 2) | Out-Null
$d.Content.Find.Execute("74×72=5328", $true, $false, $false, $false, $false, $true, 1, $false, "44×52=2288", 2) | Out-Null
$d.Content.Find.Execute("65×54=3510", $true, $false, $false, $false, $false, $true, 1, $false, "57×75=4275", 2) | Out-Null
$d.Content.Find.Execute("16×81=1296", $true, $false, $false, $false, $false, $true, 1, $false, "39×52=2028", 2) | Out-Null
$d.Content.Find.Execute("98×79=7742", $true, $false, $false, $false, $false, $true, 1, $false, "36×48=1728", 2) | Out-Null
$d.Content.Find.Execute("81×67=5427", $true, $false, $false, $false, $false, $true, 1, $false, "19×44=836", 2) | Out-Null
$d.Content.Find.Execute("81×94=7614", $true, $false, $false, $false, $false, $true, 1, $false, "94×86=8084", 2) | Out-Null
$d.Content.Find.Execute("11×32=352", $true, $false, $false, $false, $false, $true, 1, $false, "99×94=9306", 2) | Out-Null
$d.Content.Find.Execute("21×26=546", $true, $false, $false, $false, $false, $true, 1, $false, "92×40=3680", 2) | Out-Null
$d.Content.Find.Execute("99×68=6732", $true, $false, $false, $false, $false, $true, 1, $false, "56×23=1288", 2) | Out-Null
$d.Content.Find.Execute("85×88=7480", $true, $false, $false, $false, $false, $true, 1, $false, "24×90=2160", 2) | Out-Null
$d.Content.Find.Execute("64×85=5440", $true, $false, $false, $false, $false, $true, 1, $false, "33×11=363", 2) | Out-Null
